{"js": "// The Bibliografia paragraph contains a single long run of text with all\n// references concatenated together. This change breaks that text up so\n// each reference starts on its own line, using manual line breaks\n// (<w:br/>) inside the very same run (the references stay in one\n// paragraph / one run, only gaining <w:br/> separators).\nconst segments = [\n  \"- Bruice, Paula Yurkanis \u2013 Qu\u00edmica Org\u00e2nica \u2013 PEARSON Prentice Hall \u2013 S\u00e3o Paulo, 2006. Solomons, T.W.G.- Qu\u00edmica Org\u00e2nica. Volumes 1 e 2, Rio de Janeiro, Livros T\u00e9cnicos e Cient\u00edficos, 2012.\",\n  \"- Rene P. Schwarzenbach, Philip M. Gschwend, Dieter M. Imboden - Environmental Organic Chemistry \u2013 Wiley Interscience, 2016.\",\n  \"- Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14a Edi\u00e7\u00e3o, 2016. \",\n  \"- Nelson, D.; Cox, M. Princ\u00edpios de Bioqu\u00edmica de Lehninger. Artmed Editora. 6a  Edi\u00e7\u00e3o, 2014.\",\n  \"- Pratt, C.; Cornely, K. Bioqu\u00edmica essencial. Guanabara Koogan. 1a  Edi\u00e7\u00e3o, 2006. \",\n  \"- Wasserman, S.A.; Minorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora. 8 a  Edi\u00e7\u00e3o. 2010.\",\n  \"- Cooper, G.M. A C\u00e9lula \u2013 Uma Abordagem molecular. Artmed Editora Ltda. 3a  Edi\u00e7\u00e3o. 2007.\",\n  \"- Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.- Maier, R. Environmental Microbiology. Academic Press. 2000. \",\n  \"- Jordening, H.; Winter, J. Environmental Biotechnology. Concepts and Applications. Wiley-VCH. 2005. \",\n  \"- Brock, T. D. ; Madigan, M.T.; Martinko, J.M.; Dunlap, P.; Clark, D. Biology of Microorganisms. Pearson Education.12a  Edi\u00e7\u00e3o. 2009.\",\n  \"- Tortora, G.; Burdell, B.; Case, C. Microbiology. An Introduction. Pearson Benjamin Cummings. 10a  Edi\u00e7\u00e3o. 2010.\",\n];\n\nfunction xmlEscape(s) {\n  return s.replace(/&/g, \"&amp;\").replace(/</g, \"&lt;\").replace(/>/g, \"&gt;\");\n}\n\n// Build the run's inner XML: <w:t>seg0</w:t><w:br/><w:t>seg1</w:t>...\n// Any segment with leading/trailing whitespace needs xml:space=\"preserve\"\n// so the space survives round-tripping.\nlet runInner = \"\";\nfor (let i = 0; i < segments.length; i++) {\n  const seg = segments[i];\n  const preserve = seg !== seg.trim() ? ' xml:space=\"preserve\"' : \"\";\n  runInner += `<w:t${preserve}>${xmlEscape(seg)}</w:t>`;\n  if (i !== segments.length - 1) {\n    runInner += \"<w:br/>\";\n  }\n}\n\nconst flatOpcXml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  `<w:body><w:p><w:r>${runInner}</w:r></w:p></w:body>` +\n  \"</w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the bibliography paragraph: the one starting with the first\n// reference entry (Bruice).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"- Bruice, Paula Yurkanis\") === 0) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Bibliografia paragraph not found\");\n}\n\n// Replace the whole paragraph's content (a single run) with a single run\n// that contains all the references separated by manual line breaks.\nconst range = target.getRange();\nrange.insertOoxml(flatOpcXml, \"Replace\");\nawait context.sync();\n", "ps1": "# The Bibliografia paragraph holds a single run of text with every\n# reference concatenated together. This rewrites that run so each\n# reference starts on its own line, using manual line breaks (<w:br/>)\n# inside the very same run (Shift+Enter style) rather than starting new\n# paragraphs.\n\n$segments = @(\n    \"- Bruice, Paula Yurkanis \u2013 Qu\u00edmica Org\u00e2nica \u2013 PEARSON Prentice Hall \u2013 S\u00e3o Paulo, 2006. Solomons, T.W.G.- Qu\u00edmica Org\u00e2nica. Volumes 1 e 2, Rio de Janeiro, Livros T\u00e9cnicos e Cient\u00edficos, 2012.\",\n    \"- Rene P. Schwarzenbach, Philip M. Gschwend, Dieter M. Imboden - Environmental Organic Chemistry \u2013 Wiley Interscience, 2016.\",\n    \"- Madigan, M.T.; Martinko, J.M.; Bender, K.S.; Buckley, D.H.; Stahl, D.A. Microbiologia de Brock. Editora Artmed, 14a Edi\u00e7\u00e3o, 2016. \",\n    \"- Nelson, D.; Cox, M. Princ\u00edpios de Bioqu\u00edmica de Lehninger. Artmed Editora. 6a  Edi\u00e7\u00e3o, 2014.\",\n    \"- Pratt, C.; Cornely, K. Bioqu\u00edmica essencial. Guanabara Koogan. 1a  Edi\u00e7\u00e3o, 2006. \",\n    \"- Wasserman, S.A.; Minorsky, P.V.; Jackson, R.; Reece, J.; Cain, M.; Urry, L. Biologia de Campbell. Artmed Editora. 8 a  Edi\u00e7\u00e3o. 2010.\",\n    \"- Cooper, G.M. A C\u00e9lula \u2013 Uma Abordagem molecular. Artmed Editora Ltda. 3a  Edi\u00e7\u00e3o. 2007.\",\n    \"- Raven, P.H.; Evert, S.E. Biologia vegetal. Editora Guanabara Koogan, 2007.- Maier, R. Environmental Microbiology. Academic Press. 2000. \",\n    \"- Jordening, H.; Winter, J. Environmental Biotechnology. Concepts and Applications. Wiley-VCH. 2005. \",\n    \"- Brock, T. D. ; Madigan, M.T.; Martinko, J.M.; Dunlap, P.; Clark, D. Biology of Microorganisms. Pearson Education.12a  Edi\u00e7\u00e3o. 2009.\",\n    \"- Tortora, G.; Burdell, B.; Case, C. Microbiology. An Introduction. Pearson Benjamin Cummings. 10a  Edi\u00e7\u00e3o. 2010.\"\n)\n\nfunction XmlEscape($s) {\n    return $s.Replace(\"&\", \"&amp;\").Replace(\"<\", \"&lt;\").Replace(\">\", \"&gt;\")\n}\n\n# Build the run's inner XML: <w:t>seg0</w:t><w:br/><w:t>seg1</w:t>...\n# Any segment with leading/trailing whitespace needs xml:space=\"preserve\"\n# so the space survives round-tripping.\n$runInner = \"\"\nfor ($i = 0; $i -lt $segments.Count; $i++) {\n    $seg = $segments[$i]\n    $preserve = \"\"\n    if ($seg -ne $seg.Trim()) {\n        $preserve = ' xml:space=\"preserve\"'\n    }\n    $runInner += \"<w:t\" + $preserve + \">\" + (XmlEscape $seg) + \"</w:t>\"\n    if ($i -ne $segments.Count - 1) {\n        $runInner += \"<w:br/>\"\n    }\n}\n\n$flatOpcXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p><w:r>' + $runInner + '</w:r></w:p></w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n\n$d = $word.ActiveDocument\n\n# Locate the bibliography paragraph (the one starting with the first\n# reference entry).\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.StartsWith(\"- Bruice, Paula Yurkanis\")) {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Bibliografia paragraph not found\"\n}\n\n# Address the paragraph's content *excluding* its trailing paragraph mark,\n# so InsertXML replaces the run in place instead of inserting a sibling\n# paragraph.\n$pStart = $target.Range.Start\n$pText = $target.Range.Text\n$contentLen = $pText.TrimEnd([char]13, [char]7).Length\n$contentRange = $d.Range($pStart, $pStart + $contentLen)\n$contentRange.InsertXML($flatOpcXml)\n"}
